# Refined metadata to be additional tab
#
# 1. Add a new "metadata" worksheet after the existing "data" sheet,
#    containing one summary row describing the panel query that produced
#    the "data" sheet.
# 2. Refresh the per-gene "time_taken" timestamps on the "data" sheet to
#    the time of this (re-)run.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- add the "metadata" sheet right after "data" ---------------------------
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Infantile nystagmus"
$newSheet.Range("C2").Value = 246

# data_version looks numeric ("1.3") but must stay textual, like the
# geneConfidence column on the "data" sheet - force text before assigning.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.3"

$newSheet.Range("E2").Value = "2019-06-20T15:12:30.313760Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:05.507001"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/246/?format=json"

# Match the look of the "data" sheet's header/index styling (bold, centered,
# thin border) by copying the existing formats instead of inventing new
# style entries.
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# D2's forced "@" number format leaves a stray style behind - reset it back
# to the plain (unstyled) look shared by the rest of the data row by
# re-pasting the format from an unstyled cell.
$dataSheet.Range("B2").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)

# --- refresh the query timestamps on the "data" sheet -----------------------
$dataSheet.Range("F2").Value = "2021-10-05 14:21:05.510738"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:05.510746"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:05.510749"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:05.510752"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:05.510755"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:05.510758"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:05.510760"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:05.510763"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:05.510766"
$dataSheet.Range("F11").Value = "2021-10-05 14:21:05.510768"
$dataSheet.Range("F12").Value = "2021-10-05 14:21:05.510771"
$dataSheet.Range("F13").Value = "2021-10-05 14:21:05.510773"
$dataSheet.Range("F14").Value = "2021-10-05 14:21:05.510776"
$dataSheet.Range("F15").Value = "2021-10-05 14:21:05.510778"
$dataSheet.Range("F16").Value = "2021-10-05 14:21:05.510781"
$dataSheet.Range("F17").Value = "2021-10-05 14:21:05.510783"
$dataSheet.Range("F18").Value = "2021-10-05 14:21:05.510786"
$dataSheet.Range("F19").Value = "2021-10-05 14:21:05.510789"
